$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Record", "Balanço Geral", "Social", "2025-04-02T13:11", "Positivo", "Primeiro dia do mutirão de atualização do Cadúnico é marcado por tumulto. Repórter *ao vivo*. *nota da prefeitura*"),
    @("Record", "Balanço Geral", "PROCON", "2025-04-02T13:18", "Positivo", "Notícias em destaque na região. Procon Campos divulgou pesquisa de preços de produtos mais consumidos na praça. 50 itens pesquisados em 4 supermercados. "),
    @("Record", "Balanço Geral", "Social", "2025-04-02T13:17", "Neutro", "Notícias em destaque na região. Bebê encontrado em uma sacola no Parque Santa helena segue internado no Hospital Plantadores de Cana. Quando tiver alta, vai para um abrigo assistido pela prefeitura.")
)

$startRow = 7
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $data[$i][$col - 1]
    }
}
